# changes done in TC37 data and added wait
$wb = $excel.ActiveWorkbook

# Select a range on tc048 (the sheet that will lose focus) so its
# saved selection matches what gets persisted once it is no longer active.
$tc048 = $wb.Worksheets.Item("tc048")
$tc048.Activate()
$tc048.Range("A1:C2").Select()

# Insert the new "tc037" worksheet right after "tc048".
$newSheet = $wb.Worksheets.Add($null, $tc048)
$newSheet.Name = "tc037"

# Populate the header row and the data row.
$newSheet.Range("A1").Value = "Epic"
$newSheet.Range("B1").Value = "Feature"
$newSheet.Range("C1").Value = "Requirement"
$newSheet.Range("A2").Value = "Epic Mohit"
$newSheet.Range("B2").Value = "Mohit Feature"
$newSheet.Range("C2").Value = "RQ-489"

# Match formatting used by the sibling test-case tabs: wrapped header
# style with a taller row height.
$newSheet.Range("A1:C2").WrapText = $true
$newSheet.Rows.Item(1).RowHeight = 29
$newSheet.Rows.Item(2).RowHeight = 29

$newSheet.Activate()
$newSheet.Range("E16").Select()
